$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '42.482.46'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.547.21'
$ws.Range("E3").Value = '  +0.81%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '312.60'
$ws.Range("E5").Value = '  -1.18%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '100.68'
$ws.Range("E6").Value = '  +4.96%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.569'
$ws.Range("E7").Value = '  -0.77%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.530'
$ws.Range("E9").Value = '  -1.18%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.06'
$ws.Range("E10").Value = '  +1.83%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0804'
$ws.Range("E11").Value = '  -0.60%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.41'
$ws.Range("E12").Value = '  -0.60%  '
$ws.Range("E13").Value = '  -0.11%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.937.76'
$ws.Range("E14").Value = '  +0.62%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.98'
$ws.Range("E15").Value = '  +6.82%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.569.05'
$ws.Range("E16").Value = '  +0.23%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.841'
$ws.Range("E17").Value = '  -0.37%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '42.526.71'
$ws.Range("E18").Value = '  -0.76%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.80'
$ws.Range("E19").Value = '  +0.13%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0953'
$ws.Range("E20").Value = '  -0.74%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.32'
$ws.Range("E21").Value = '  -1.34%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '69.14'
$ws.Range("E22").Value = '  -0.27%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '244.53'
$ws.Range("E23").Value = '  -3.29%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.93'
$ws.Range("E24").Value = '  -0.96%  '
$ws.Range("E25").Value = '  +1.09%  '
$ws.Range("B26").Value = 'EthereumClassic'
$ws.Range("C26").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '26.39'
$ws.Range("E26").Value = '  -0.93%  '
$ws.Range("B27").Value = 'Dai'
$ws.Range("C27").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.00'
$ws.Range("E27").Value = '  -0.02%  '
$ws.Range("B28").Value = 'Toncoin'
$ws.Range("C28").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.36'
$ws.Range("E28").Value = '  -1.66%  '
$ws.Range("B29").Value = 'InjectiveProtocol'
$ws.Range("C29").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '40.34'
$ws.Range("E29").Value = '  -0.06%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '10.12'
$ws.Range("E30").Value = '  -0.95%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '158.19'
$ws.Range("E31").Value = '  +1.24%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.70'
$ws.Range("E32").Value = '  -2.67%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.78'
$ws.Range("E33").Value = '  +14.37%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0804'
$ws.Range("E34").Value = '  +1.38%  '
$ws.Range("B35").Value = 'ARBITRUM'
$ws.Range("C35").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.06'
$ws.Range("E35").Value = '  -0.86%  '
$ws.Range("B36").Value = 'WEMIXToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.63'
$ws.Range("E36").Value = '  -2.93%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.20'
$ws.Range("E37").Value = '  -3.51%  '
$ws.Range("E38").Value = '  -4.45%  '
$ws.Range("E39").Value = '  -0.69%  '
$ws.Range("E40").Value = '  +0.04%  '
$ws.Range("E41").Value = '  +9.91%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '22.16'
$ws.Range("E42").Value = '  +2.15%  '
$ws.Range("B43").Value = 'NEARProtocol'
$ws.Range("C43").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.33'
$ws.Range("E43").Value = '  +2.50%  '
$ws.Range("B44").Value = 'FirstDigitalUSD'
$ws.Range("C44").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.00'
$ws.Range("E44").Value = '  +0.07%  '
$ws.Range("E45").Value = '  -1.45%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.966.54'
$ws.Range("E46").Value = '  -1.25%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.88'
$ws.Range("E47").Value = '  -1.81%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.789.46'
$ws.Range("E48").Value = '  +0.57%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '81.13'
$ws.Range("E49").Value = '  -3.68%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.193'
$ws.Range("E50").Value = '  +1.29%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '73.17'
$ws.Range("E51").Value = '  -1.31%  '
